$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.549.38'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '3.837.37'
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'600.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = "'163.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("D7").Value = '3.836.77'
$ws.Range("E7").Value = '  +2.57%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").Value = "'6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("E13").Value = '  -3.44%  '
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = '4.483.92'
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = '3.855.62'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("D17").Value = '68.746.95'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = "'7.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").Value = "'11.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = "'486.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.16%  '
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Value = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").Value = "'2.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").Value = '3.990.74'
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").Value = "'2.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.30%  '
$ws.Range("D34").Value = "'31.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").Value = '3.783.78'
$ws.Range("E35").Value = '  +3.07%  '
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = "'0.318"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("D43").Value = "'431.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").Value = "'48.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").Value = '2.844.60'
$ws.Range("E48").Value = '  +2.15%  '
$ws.Range("D49").Value = "'142.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").Value = "'0.0357"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = "'25.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +13.62%  '
